$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 19898.5
$ws.Range("J3").Value = 19898.5
$ws.Range("L3").Value = 19898.5
$ws.Range("N3").Value = -20126.5
$ws.Range("H63").Value = 0.0
$ws.Range("J63").Value = 0.0
$ws.Range("L63").Value = 0.0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 10199.08
$ws.Range("I64").Value = 6908.0
$ws.Range("J64").Value = 12784.929
$ws.Range("K64").Value = 6908.0
$ws.Range("L64").Value = 12784.929
$ws.Range("M64").Value = -6660.0
$ws.Range("N64").Value = -13280.929
$ws.Range("H66").Value = 0.0
$ws.Range("J66").Value = 0.0
$ws.Range("L66").Value = 0.0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 10199.08
$ws.Range("I67").Value = 6908.0
$ws.Range("J67").Value = 12784.929
$ws.Range("K67").Value = 6908.0
$ws.Range("L67").Value = 12784.929
$ws.Range("M67").Value = -6050.0
$ws.Range("N67").Value = -14500.929
$ws.Range("H80").Value = 1802.6154
$ws.Range("I80").Value = 787.1111
$ws.Range("J80").Value = 4087.5
$ws.Range("K80").Value = 2361.3333
$ws.Range("L80").Value = 12262.5
$ws.Range("M80").Value = -1363.3333
$ws.Range("N80").Value = -14258.5
$ws.Range("H81").Value = 75000.0
$ws.Range("J81").Value = 75000.0
$ws.Range("L81").Value = 75000.0
$ws.Range("N81").Value = -76996.0
$ws.Range("H82").Value = 468.5
$ws.Range("I82").Value = 468.5
$ws.Range("J82").Value = 0.0
$ws.Range("K82").Value = 1405.5
$ws.Range("L82").Value = 0.0
$ws.Range("M82").Value = -999.5
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 1802.6154
$ws.Range("I83").Value = 787.1111
$ws.Range("J83").Value = 4087.5
$ws.Range("K83").Value = 7083.9999
$ws.Range("L83").Value = 36787.5
$ws.Range("M83").Value = -2091.9999
$ws.Range("N83").Value = -46771.5
$ws.Range("H84").Value = 75000.0
$ws.Range("J84").Value = 75000.0
$ws.Range("L84").Value = 225000.0
$ws.Range("N84").Value = -234984.0
$ws.Range("H85").Value = 468.5
$ws.Range("I85").Value = 468.5
$ws.Range("J85").Value = 0.0
$ws.Range("K85").Value = 1405.5
$ws.Range("L85").Value = 0.0
$ws.Range("M85").Value = -1.5
$ws.Range("N85").ClearContents()
$ws.Range("H92").Value = 505.8889
$ws.Range("I92").Value = 468.5
$ws.Range("J92").Value = 805.0
$ws.Range("K92").Value = 468.5
$ws.Range("L92").Value = 805.0
$ws.Range("M92").Value = 779.5
$ws.Range("N92").Value = -3301.0
$ws.Range("H93").Value = 50000.0
$ws.Range("J93").Value = 50000.0
$ws.Range("L93").Value = 50000.0
$ws.Range("N93").Value = -54992.0
$ws.Range("H101").Value = 2183.3333
$ws.Range("J101").Value = 2000.0
$ws.Range("L101").Value = 6000.0
$ws.Range("N101").Value = -9244.0
$ws.Range("H102").Value = 19898.5
$ws.Range("J102").Value = 19898.5
$ws.Range("L102").Value = 19898.5
$ws.Range("N102").Value = -26388.5
$ws.Range("H103").Value = 3285.7144
$ws.Range("J103").Value = 5000.0
$ws.Range("L103").Value = 15000.0
$ws.Range("N103").Value = -16172.0
$ws.Range("H105").Value = 14998.5
$ws.Range("J105").Value = 14998.5
$ws.Range("L105").Value = 14998.5
$ws.Range("N105").Value = -21986.5
$ws.Range("H106").Value = 0.0
$ws.Range("I106").Value = 0.0
$ws.Range("K106").Value = 0.0
$ws.Range("M106").ClearContents()
$ws.Range("H107").Value = 376.15384
$ws.Range("I107").Value = 328.75
$ws.Range("J107").Value = 945.0
$ws.Range("K107").Value = 328.75
$ws.Range("L107").Value = 945.0
$ws.Range("M107").Value = 1591.25
$ws.Range("N107").Value = -4785.0
$ws.Range("H109").Value = 0.0
$ws.Range("J109").Value = 0.0
$ws.Range("L109").Value = 0.0
$ws.Range("N109").ClearContents()
$ws.Range("H138").Value = 5344.773
$ws.Range("J138").Value = 5599.5835
$ws.Range("L138").Value = 16798.7505
$ws.Range("N138").Value = -27078.7505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 3343815.2
$ws.Range("J8").Value = 4295.3335
$ws.Range("L8").Value = 4295.3335
$ws.Range("N8").Value = -4583.3335
$ws.Range("H30").Value = 2139.0
$ws.Range("J30").Value = 4269.0
$ws.Range("L30").Value = 4269.0
$ws.Range("N30").Value = -4569.0
$ws.Range("H63").Value = 2150.4
$ws.Range("I63").Value = 1445.0
$ws.Range("K63").Value = 1445.0
$ws.Range("M63").Value = -759.0
$ws.Range("H66").Value = 2150.4
$ws.Range("I66").Value = 1445.0
$ws.Range("K66").Value = 7225.0
$ws.Range("M66").Value = -3793.0
$ws.Range("H74").Value = 1471.7858
$ws.Range("I74").Value = 1416.3636
$ws.Range("K74").Value = 1416.3636
$ws.Range("M74").Value = -542.3635999999999
$ws.Range("H77").Value = 1471.7858
$ws.Range("I77").Value = 1416.3636
$ws.Range("K77").Value = 7081.817999999999
$ws.Range("M77").Value = -2713.817999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1715.0
$ws.Range("I5").Value = 72.5
$ws.Range("J5").Value = 5000.0
$ws.Range("K5").Value = 72.5
$ws.Range("L5").Value = 5000.0
$ws.Range("M5").Value = 40.5
$ws.Range("N5").Value = -5226.0
$ws.Range("H19").Value = 5000.0
$ws.Range("I19").Value = 5000.0
$ws.Range("K19").Value = 5000.0
$ws.Range("M19").Value = -4827.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 4405.75
$ws.Range("I10").Value = 2187.25
$ws.Range("J10").Value = 6624.25
$ws.Range("K10").Value = 2187.25
$ws.Range("L10").Value = 6624.25
$ws.Range("M10").Value = -2048.25
$ws.Range("N10").Value = -6902.25
$ws.Range("H31").Value = 2994.25
$ws.Range("I31").Value = 2994.25
$ws.Range("K31").Value = 2994.25
$ws.Range("M31").Value = -2699.25
$ws.Range("H34").Value = 2994.25
$ws.Range("I34").Value = 2994.25
$ws.Range("K34").Value = 2994.25
$ws.Range("M34").Value = -2792.25
$ws.Range("H132").Value = 6077.9
$ws.Range("I132").Value = 3873.125
$ws.Range("K132").Value = 11619.375
$ws.Range("M132").Value = -9089.375
$ws.Range("H134").Value = 2236.3333
$ws.Range("I134").Value = 1719.25
$ws.Range("K134").Value = 5157.75
$ws.Range("M134").Value = -2622.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 811.625
$ws.Range("I4").Value = 856.2857
$ws.Range("J4").Value = 499.0
$ws.Range("K4").Value = 2568.8571
$ws.Range("L4").Value = 1497.0
$ws.Range("M4").Value = -2456.8571
$ws.Range("N4").Value = -1721.0
$ws.Range("H8").Value = 312.57144
$ws.Range("I8").Value = 312.57144
$ws.Range("K8").Value = 937.71432
$ws.Range("M8").Value = -798.71432
$ws.Range("H55").Value = 882.8333
$ws.Range("J55").Value = 1124.375
$ws.Range("L55").Value = 3373.125
$ws.Range("N55").Value = -3727.125
$ws.Range("H128").Value = 499991.0
$ws.Range("I128").Value = 499991.0
$ws.Range("K128").Value = 1499973.0
$ws.Range("M128").Value = -1494993.0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1714.0
$ws.Range("I11").Value = 2141.5
$ws.Range("K11").Value = 2141.5
$ws.Range("M11").Value = -2002.5
$ws.Range("H21").Value = 5000.0
$ws.Range("I21").Value = 5000.0
$ws.Range("K21").Value = 5000.0
$ws.Range("M21").Value = -4827.0
$ws.Range("H30").Value = 5000.0
$ws.Range("I30").Value = 5000.0
$ws.Range("K30").Value = 5000.0
$ws.Range("M30").Value = -4895.0
$ws.Range("H70").Value = 10649.111
$ws.Range("I70").Value = 19949.334
$ws.Range("J70").Value = 5999.0
$ws.Range("K70").Value = 19949.334
$ws.Range("L70").Value = 5999.0
$ws.Range("M70").Value = -19679.334
$ws.Range("N70").Value = -6539.0
$ws.Range("H73").Value = 10649.111
$ws.Range("I73").Value = 19949.334
$ws.Range("J73").Value = 5999.0
$ws.Range("K73").Value = 19949.334
$ws.Range("L73").Value = 5999.0
$ws.Range("M73").Value = -19013.334
$ws.Range("N73").Value = -7871.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 930.9231
$ws.Range("J55").Value = 930.9231
$ws.Range("L55").Value = 930.9231
$ws.Range("N55").Value = -1276.9231
$ws.Range("H109").Value = 26000.0
$ws.Range("I109").Value = 14000.0
$ws.Range("J109").Value = 38000.0
$ws.Range("K109").Value = 14000.0
$ws.Range("L109").Value = 38000.0
$ws.Range("M109").Value = -12613.0
$ws.Range("N109").Value = -40774.0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 19999.0
$ws.Range("I2").Value = 19999.0
$ws.Range("K2").Value = 19999.0
$ws.Range("M2").Value = -19887.0
